# This script applies the data-refresh changes described by the commit
# "Update gh-pages to output generated at 456a3b4" to the 北京-漫展信息
# workbook. It updates the "想去人数" (want-to-go count) and "最低票价"
# (minimum ticket price) columns, a couple of addresses, and two cover
# image URLs across the 展览 and 全部类型 sheets, plus a single count on
# the 演出 sheet.
#
# Note: this runtime's parameter binder does not reliably bind named
# (-Param value) arguments, so helper functions below use positional
# parameters only.

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($SheetName, $CellRef, $Value)
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $Value
}

# ---------------------------------------------------------------------
# Sheet: 展览
# ---------------------------------------------------------------------
$sheet = "展览"

Set-CellValue $sheet "F4"  2940
Set-CellValue $sheet "F7"  2298
Set-CellValue $sheet "F8"  1607
Set-CellValue $sheet "F11" 103
Set-CellValue $sheet "F12" 2636
Set-CellValue $sheet "F14" 1480
Set-CellValue $sheet "F15" 6928
Set-CellValue $sheet "F17" 7082
Set-CellValue $sheet "F19" 2870
Set-CellValue $sheet "F21" 3443
Set-CellValue $sheet "F22" 217
Set-CellValue $sheet "F23" 147
Set-CellValue $sheet "I23" "//i0.hdslb.com/bfs/openplatform/202405/mwy8WpFC1715567332261.jpeg"
Set-CellValue $sheet "F24" 1820
Set-CellValue $sheet "F25" 74
Set-CellValue $sheet "I25" "//i1.hdslb.com/bfs/openplatform/202405/gVA8JEt21715567367779.jpeg"
Set-CellValue $sheet "F26" 289
Set-CellValue $sheet "F29" 171
Set-CellValue $sheet "F31" 383
Set-CellValue $sheet "G31" 70
Set-CellValue $sheet "F33" 2526
Set-CellValue $sheet "D35" "金蝉西路甲1号 北京酷车国际汇展中心"
Set-CellValue $sheet "F35" 158
Set-CellValue $sheet "F36" 366
Set-CellValue $sheet "F37" 1007
Set-CellValue $sheet "F38" 203
Set-CellValue $sheet "F39" 456
Set-CellValue $sheet "F40" 505
Set-CellValue $sheet "G40" 75

# ---------------------------------------------------------------------
# Sheet: 演出
# ---------------------------------------------------------------------
$sheet = "演出"

Set-CellValue $sheet "F10" 4

# ---------------------------------------------------------------------
# Sheet: 全部类型
# ---------------------------------------------------------------------
$sheet = "全部类型"

Set-CellValue $sheet "F6"  2940
Set-CellValue $sheet "F8"  2298
Set-CellValue $sheet "F9"  1607
Set-CellValue $sheet "F12" 103
Set-CellValue $sheet "F14" 2636
Set-CellValue $sheet "F15" 1480
Set-CellValue $sheet "F18" 4
Set-CellValue $sheet "F20" 6928
Set-CellValue $sheet "F22" 7082
Set-CellValue $sheet "F24" 2870
Set-CellValue $sheet "F26" 3443
Set-CellValue $sheet "F28" 217
Set-CellValue $sheet "F31" 1820
Set-CellValue $sheet "F34" 289
Set-CellValue $sheet "F37" 171
Set-CellValue $sheet "F39" 383
Set-CellValue $sheet "G39" 70
Set-CellValue $sheet "F41" 2526
Set-CellValue $sheet "D43" "金蝉西路甲1号 北京酷车国际汇展中心"
Set-CellValue $sheet "F43" 158
Set-CellValue $sheet "F45" 366
Set-CellValue $sheet "F46" 1007
Set-CellValue $sheet "F47" 203
Set-CellValue $sheet "F48" 456
Set-CellValue $sheet "F49" 505
Set-CellValue $sheet "G49" 75

$wb.Save()
